$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column header "OSMO_DEF" in F1, matching the style of the other headers (B1:E1)
$ws.Range("F1").Value = "OSMO_DEF"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Add new values in column F for data rows
$ws.Range("F2").Value = "[]"
$ws.Range("F3").Value = "[]"
